$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B15").Value = "James"
$ws.Range("C15").Value = "ResultsProcessing?"
$ws.Range("B24").Value = "James"

$ws.Range("C15").Font.Name = "Consolas"
$ws.Range("C15").Font.Size = 9.8
$ws.Range("C15").Font.Color = 0
$ws.Range("C15").VerticalAlignment = -4108

$ws.Range("C8").Select()
